# 自动更新Excel文件 - apply daily countdown update to 剩余(remaining) / 开始时间(start date)
# columns, and append newly discovered shops at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the new value for column E ("剩余", remaining) for the given row,
# and (when present) the new value for column F ("开始时间", start date) when the
# item was restocked (remaining hit 0 and got refilled back up to the total
# day count, with the start date reset to the most recent count date).
$updates = @(
    @{Row=2; E=7},
    @{Row=3; E=7},
    @{Row=4; E=7},
    @{Row=5; E=3},
    @{Row=6; E=7},
    @{Row=7; E=3},
    @{Row=8; E=7},
    @{Row=9; E=3},
    @{Row=10; E=7; F=20251013},
    @{Row=11; E=7},
    @{Row=12; E=3},
    @{Row=13; E=7},
    @{Row=14; E=7},
    @{Row=15; E=7},
    @{Row=16; E=7},
    @{Row=17; E=3},
    @{Row=18; E=6},
    @{Row=19; E=6},
    @{Row=20; E=6},
    @{Row=21; E=6},
    @{Row=22; E=3},
    @{Row=23; E=3},
    @{Row=24; E=3},
    @{Row=25; E=3},
    @{Row=26; E=3},
    @{Row=27; E=1},
    @{Row=28; E=6},
    @{Row=29; E=6},
    @{Row=30; E=6},
    @{Row=31; E=6},
    @{Row=32; E=6},
    @{Row=33; E=6},
    @{Row=34; E=6},
    @{Row=35; E=6},
    @{Row=37; E=6},
    @{Row=38; E=6},
    @{Row=39; E=6},
    @{Row=40; E=7; F=20251013},
    @{Row=41; E=7; F=20251013},
    @{Row=42; E=6},
    @{Row=43; E=3},
    @{Row=44; E=7; F=20251013},
    @{Row=45; E=3},
    @{Row=46; E=7; F=20251013},
    @{Row=47; E=6},
    @{Row=48; E=7; F=20251013},
    @{Row=49; E=1},
    @{Row=50; E=1},
    @{Row=51; E=1},
    @{Row=52; E=1},
    @{Row=53; E=1},
    @{Row=54; E=1},
    @{Row=55; E=1},
    @{Row=56; E=1},
    @{Row=57; E=1},
    @{Row=58; E=5},
    @{Row=59; E=5},
    @{Row=60; E=5},
    @{Row=61; E=1},
    @{Row=62; E=5},
    @{Row=63; E=5},
    @{Row=64; E=5},
    @{Row=65; E=6},
    @{Row=66; E=6},
    @{Row=67; E=6},
    @{Row=68; E=6},
    @{Row=69; E=6},
    @{Row=70; E=7},
    @{Row=71; E=7},
    @{Row=72; E=7},
    @{Row=73; E=7},
    @{Row=74; E=7},
    @{Row=75; E=7},
    @{Row=76; E=7},
    @{Row=77; E=10; F=20251013},
    @{Row=78; E=10; F=20251013},
    @{Row=79; E=10; F=20251013},
    @{Row=80; E=10; F=20251013},
    @{Row=81; E=10; F=20251013},
    @{Row=82; E=10; F=20251013},
    @{Row=83; E=10; F=20251013},
    @{Row=84; E=10; F=20251013},
    @{Row=85; E=10; F=20251013},
    @{Row=86; E=10; F=20251013},
    @{Row=87; E=7; F=20251013},
    @{Row=88; E=7; F=20251013},
    @{Row=89; E=7; F=20251013},
    @{Row=90; E=7; F=20251013},
    @{Row=91; E=3},
    @{Row=92; E=7; F=20251013},
    @{Row=93; E=10; F=20251013},
    @{Row=94; E=3},
    @{Row=95; E=9}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 5).Value = $u.E
    if ($u.ContainsKey("F")) {
        $ws.Cells.Item($u.Row, 6).Value = $u.F
    }
}

# Append the newly scraped shops as rows 96-99.
$newRows = @(
    @{Row=96; A=95; B="俏大姐";   C="唐王河路（小湖美食城）"; D=10; E=7; F=20251010; G="";       H="小桶1个"; I=""},
    @{Row=97; A=96; B="西子居";   C="体育馆";                 D=10; E=7; F=20251010; G="大桶2个"; H="";       I=""},
    @{Row=98; A=97; B="德胜园";   C="邾国大道";               D=10; E=7; F=20251010; G="大桶1个"; H="";       I=""},
    @{Row=99; A=98; B="老鲁味";   C="唐王河路（小湖美食城）"; D=10; E=7; F=20251010; G="大桶1个"; H="";       I=""}
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    if ($r.G -ne "") { $ws.Cells.Item($r.Row, 7).Value = $r.G }
    if ($r.H -ne "") { $ws.Cells.Item($r.Row, 8).Value = $r.H }
    if ($r.I -ne "") { $ws.Cells.Item($r.Row, 9).Value = $r.I }
}
